$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert a new "Min Voltage" column before the existing
# Price/Quantity/Total columns (F,G,H -> G,H,I), shifting values+style right
# one cell and writing the new header into F1.
$ws.Range("I1").Value2 = $ws.Range("H1").Value2
$ws.Range("I1").Font.Size = $ws.Range("H1").Font.Size
$ws.Range("I1").Font.Color = $ws.Range("H1").Font.Color

$ws.Range("H1").Value2 = $ws.Range("G1").Value2
$ws.Range("G1").Value2 = $ws.Range("F1").Value2
$ws.Range("F1").Value2 = "Min Voltage"

# New column F is a bit narrower than the others.
$ws.Columns(6).ColumnWidth = 12.3

# --- Rows 6-9 (the various electrolytic-cap footprints that are no longer
# in play) are now hidden rather than deleted.
$ws.Rows(6).Hidden = $true
$ws.Rows(7).Hidden = $true
$ws.Rows(8).Hidden = $true
$ws.Rows(9).Hidden = $true

# --- Row 10 switches from the old electrolytic 22nF cap to the ceramic
# replacement, now with a part number + package.
$ws.Range("C10").Value2 = "22nf Ceramic Cap"
$ws.Range("D10").Value2 = "CL10B223KB8NNNC (JLC)"
$ws.Range("E10").Value2 = "0603(imperial)"

# --- Row 11 (new): 10uF ceramic cap.
$ws.Range("B11").Value2 = "C8,C9,C12,C13,C22,C54,C57"
$ws.Range("C11").Value2 = "10uF Ceramic Cap"
$ws.Range("D11").Value2 = "GRM21BR61H106KE43L (at JLC)"
$ws.Range("E11").Value2 = "0805(imperial)"

# --- Row 12 (new): 2.2uF ceramic cap.
$ws.Range("B12").Value2 = "C16"
$ws.Range("C12").Value2 = "2.2uF Ceramic Cap"
$ws.Range("D12").Value2 = "0805F225M500NT (JLC)"
$ws.Range("E12").Value2 = "0805(imperial)"

[void]$ws.Range("E12").Select()
